# Physical asset import template: add an "Asset Owner Name" column
# between "Asset Number*" (B) and "Asset Type" (existing C, becomes D).
#
# This inserts a new column at C (shifting Asset Type / Asset Value and
# the trailing formatted blank cell one column to the right, from C/D/K
# to D/E/L) and labels the new header cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the current column C ("Asset Type"),
# shifting it (and everything after it) one column to the right.
$ws.Columns("C").Insert()

# Label the newly inserted header cell.
$ws.Range("C1").Value = "Asset Owner Name"

# Match the saved selection state (active cell D5).
[void]$ws.Range("D5").Select()
